# Weekly refresh of fruit/vegetable price records (commit: "Fruta / hortaliza, semanal").
# The report rows (2-43) are resequenced: Fecha/Calidad/Volumen/Precio columns are
# updated in place per row to reflect the new weekly snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44432
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 1200
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = 1250
$ws.Range("P2").Value = 417

# Row 3
$ws.Range("D3").Value = 44432
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 950
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = 975
$ws.Range("P3").Value = 325

# Row 4
$ws.Range("D4").Value = 44428
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 1800
$ws.Range("M4").Value = 1650
$ws.Range("P4").Value = 550

# Row 5
$ws.Range("D5").Value = 44460
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 1400
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 1450
$ws.Range("P5").Value = 483

# Row 6
$ws.Range("D6").Value = 44174
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 500
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = 550
$ws.Range("P6").Value = 183

# Row 7
$ws.Range("D7").Value = 44397
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 1400
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = 1450
$ws.Range("P7").Value = 483

# Row 8
$ws.Range("D8").Value = 44475
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 250
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1200
$ws.Range("M8").Value = 1100
$ws.Range("P8").Value = 367

# Row 9
$ws.Range("D9").Value = 44431
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 1300
$ws.Range("M9").Value = 1150
$ws.Range("P9").Value = 383

# Row 10
$ws.Range("D10").Value = 44391
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 1800
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 1900
$ws.Range("P10").Value = 633

# Row 11
$ws.Range("D11").Value = 44364
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 1700
$ws.Range("L11").Value = 1800
$ws.Range("M11").Value = 1750
$ws.Range("P11").Value = 583

# Row 12
$ws.Range("D12").Value = 44364
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 1400
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = 1450
$ws.Range("P12").Value = 483

# Row 13
$ws.Range("D13").Value = 44249
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 1600
$ws.Range("M13").Value = 1550
$ws.Range("P13").Value = 517

# Row 14
$ws.Range("D14").Value = 44489
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 1400
$ws.Range("L14").Value = 1500
$ws.Range("M14").Value = 1450
$ws.Range("P14").Value = 483

# Row 15
$ws.Range("D15").Value = 44300
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 250
$ws.Range("K15").Value = 1600
$ws.Range("L15").Value = 1800
$ws.Range("M15").Value = 1700
$ws.Range("P15").Value = 567

# Row 16
$ws.Range("D16").Value = 44467
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = 1100
$ws.Range("P16").Value = 367

# Row 17
$ws.Range("D17").Value = 44447
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 1100
$ws.Range("L17").Value = 1200
$ws.Range("M17").Value = 1150
$ws.Range("P17").Value = 383

# Row 18
$ws.Range("D18").Value = 44224
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 1400
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = 1450
$ws.Range("P18").Value = 483

# Row 19
$ws.Range("D19").Value = 44224
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 160
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 1200
$ws.Range("M19").Value = 1100
$ws.Range("P19").Value = 367

# Row 20
$ws.Range("D20").Value = 44435
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 450
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 1300
$ws.Range("M20").Value = 1194
$ws.Range("P20").Value = 398

# Row 21
$ws.Range("D21").Value = 44435
$ws.Range("I21").Value = "Segunda"
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 950
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 975
$ws.Range("P21").Value = 325

# Row 22
$ws.Range("D22").Value = 44327
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 1400
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = 1450
$ws.Range("P22").Value = 483

# Row 23
$ws.Range("D23").Value = 44161
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 600
$ws.Range("L23").Value = 700
$ws.Range("M23").Value = 650
$ws.Range("P23").Value = 217

# Row 24
$ws.Range("D24").Value = 44161
$ws.Range("I24").Value = "Segunda"
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = 500
$ws.Range("L24").Value = 600
$ws.Range("M24").Value = 550
$ws.Range("P24").Value = 183

# Row 25
$ws.Range("D25").Value = 44417
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 250
$ws.Range("K25").Value = 1800
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = 1900
$ws.Range("P25").Value = 633

# Row 26
$ws.Range("D26").Value = 44417
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 200
$ws.Range("K26").Value = 1500
$ws.Range("L26").Value = 1600
$ws.Range("M26").Value = 1550
$ws.Range("P26").Value = 517

# Row 27
$ws.Range("D27").Value = 44278
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 140
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = 2250
$ws.Range("P27").Value = 750

# Row 28
$ws.Range("D28").Value = 44278
$ws.Range("I28").Value = "Segunda"
$ws.Range("J28").Value = 200
$ws.Range("K28").Value = 1500
$ws.Range("L28").Value = 1800
$ws.Range("M28").Value = 1650
$ws.Range("P28").Value = 550

# Row 29
$ws.Range("D29").Value = 44385
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 2000
$ws.Range("L29").Value = 2300
$ws.Range("M29").Value = 2150
$ws.Range("P29").Value = 717

# Row 30
$ws.Range("D30").Value = 44398
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 300
$ws.Range("K30").Value = 1700
$ws.Range("L30").Value = 1800
$ws.Range("M30").Value = 1750
$ws.Range("P30").Value = 583

# Row 31
$ws.Range("D31").Value = 44376
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 340
$ws.Range("K31").Value = 1400
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = 1471
$ws.Range("P31").Value = 490

# Row 32
$ws.Range("D32").Value = 44280
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 200
$ws.Range("K32").Value = 1800
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = 1900
$ws.Range("P32").Value = 633

# Row 33
$ws.Range("D33").Value = 44280
$ws.Range("I33").Value = "Segunda"
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 1400
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = 1450
$ws.Range("P33").Value = 483

# Row 34
$ws.Range("D34").Value = 44166
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 250
$ws.Range("K34").Value = 900
$ws.Range("L34").Value = 1000
$ws.Range("M34").Value = 950
$ws.Range("P34").Value = 317

# Row 35
$ws.Range("D35").Value = 44333
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 200
$ws.Range("K35").Value = 1500
$ws.Range("L35").Value = 1700
$ws.Range("M35").Value = 1600
$ws.Range("P35").Value = 533

# Row 36
$ws.Range("D36").Value = 44306
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 200
$ws.Range("K36").Value = 2400
$ws.Range("L36").Value = 2500
$ws.Range("M36").Value = 2450
$ws.Range("P36").Value = 817

# Row 37
$ws.Range("D37").Value = 44481
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 250
$ws.Range("K37").Value = 950
$ws.Range("L37").Value = 1000
$ws.Range("M37").Value = 975
$ws.Range("P37").Value = 325

# Row 38
$ws.Range("D38").Value = 44494
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 900
$ws.Range("L38").Value = 1000
$ws.Range("M38").Value = 950
$ws.Range("P38").Value = 317

# Row 39
$ws.Range("D39").Value = 44342
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 200
$ws.Range("K39").Value = 2000
$ws.Range("L39").Value = 2200
$ws.Range("M39").Value = 2100
$ws.Range("P39").Value = 700

# Row 42
$ws.Range("D42").Value = 44295
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 200
$ws.Range("K42").Value = 1500
$ws.Range("L42").Value = 1800
$ws.Range("M42").Value = 1650
$ws.Range("P42").Value = 550

# Row 43
$ws.Range("D43").Value = 44302
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 200
$ws.Range("K43").Value = 1400
$ws.Range("L43").Value = 1500
$ws.Range("M43").Value = 1450
$ws.Range("P43").Value = 483
